$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.249.46"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "'2.234.48"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'243.37"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").Value = "'74.49"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").Value = "'42.69"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").Value = "'0.0963"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").Value = "'6.98"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "'2.569.88"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "'14.35"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "'0.840"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "'2.231.19"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "'42.109.04"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("D21").Value = "'72.95"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "'11.18"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").Value = "'230.92"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  -5.92%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'11.39"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "'3.63"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("D30").Value = "'167.09"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'20.65"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -7.31%  "
$ws.Range("D33").Value = "'0.0805"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").Value = "'30.10"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  -7.21%  "
$ws.Range("E37").Value = "  -6.98%  "
$ws.Range("D38").Value = "'0.0305"
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").Value = "'13.27"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'2.14"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").Value = "'5.72"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'65.18"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").Value = "'0.200"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "'8.73"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").Value = "'104.72"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "'1.13"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Value = "'2.438.90"
$ws.Range("E51").Value = "  -0.70%  "
